# TradingModel - 2021/11/12 data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) of column A down to the two new rows (9 & 10)
# before touching any values, using the existing A8 cell as the format source.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Drop the old "EachCalculateProfit" (C) / "TodayClose" (D) columns entirely -
# TodayClose now becomes the new column C. Keep column C's existing formatting
# (e.g. the bold/bordered header style on C1) and only wipe column D outright
# since it disappears completely.
$ws.Range("C1:C8").ClearContents()
$ws.Range("D1:D8").Clear()

# Header row
$ws.Range("C1").Value = "TodayClose"

# Row 2 : Stock 1711
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1711
$ws.Range("C2").Value = 28

# Row 3 : Stock 2436 (new stock swapped in here)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2436
$ws.Range("C3").Value = 102.5

# Row 4 : Stock 3033
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3033
$ws.Range("C4").Value = 31.65

# Row 5 : Stock 3035
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3035
$ws.Range("C5").Value = 189

# Row 6 : Stock 3141
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 3141
$ws.Range("C6").Value = "'267.00"
$ws.Range("C6").Style = "Normal"

# Row 7 : Stock 3189
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 3189
$ws.Range("C7").Value = 251

# Row 8 : Stock 3588 (new)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 3588
$ws.Range("C8").Value = 161

# Row 9 : Stock 6104 (new)
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = 6104
$ws.Range("C9").Value = "'165.50"
$ws.Range("C9").Style = "Normal"

# Row 10 : Stock 6411
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = 6411
$ws.Range("C10").Value = "'264.00"
$ws.Range("C10").Style = "Normal"
